$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny floating-point correction on existing row 96, column A (date serial)
$ws.Cells.Item(96, 1).Value = 44409.76782184259

# Append new row 97 with newly retrieved data
$ws.Cells.Item(97, 1).Value = 44410.77463201572
$ws.Cells.Item(97, 2).Value = 79918
$ws.Cells.Item(97, 3).Value = 67637
$ws.Cells.Item(97, 4).Value = 3612
$ws.Cells.Item(97, 5).Value = 2229
$ws.Cells.Item(97, 6).Value = 1626
$ws.Cells.Item(97, 7).Value = 21116
$ws.Cells.Item(97, 8).Value = 1620
$ws.Cells.Item(97, 9).Value = 903
$ws.Cells.Item(97, 10).Value = 200
